$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$text = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.71 = 6401.37 pesos`n✅ 6401.37 pesos = 1.7 = 916.8 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $text

# --- Sheet "tasas": update the numeric values ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 583.5
$wsTasas.Range("O10").Value = 3735.2
$wsTasas.Range("N12").Value = 3760
$wsTasas.Range("O12").Value = 538.506
